$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing parameter row 12 (SM_Pos_Abstand): value 160 -> 150, highlight fill ---
$ws.Range("B12").Value = 150

# Build the "Accent1, Lighter 80%" fill on a scratch cell first and copy the
# resulting format over, so the real target range gets a single clean style
# index (avoids leaving extra intermediate cellXf entries behind).
$scratch = $ws.Range("Z1")
$scratch.Interior.ThemeColor = 5
$scratch.Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$scratch.Clear()

# --- Append new parameter rows 29-35 ("Laengsstreben" / angle-construction values) ---

# Row 29: Rundung_Platte
$ws.Range("A29").Value = "Rundung_Platte"
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = "mm"

# Row 30: GM_Schrauben_Pos_Abstand_kurz (with comment "m2,5")
$ws.Range("A30").Value = "GM_Schrauben_Pos_Abstand_kurz"
$ws.Range("B30").Value = 18.37
$ws.Range("C30").Value = "mm"

# Row 31: GM_Schrauben_Pos_Abstand_lang
$ws.Range("A31").Value = "GM_Schrauben_Pos_Abstand_lang"
$ws.Range("B31").Value = 34.37
$ws.Range("C31").Value = "mm"

# Row 32: GM_Stift_Durchmesser
$ws.Range("A32").Value = "GM_Stift_Durchmesser"
$ws.Range("B32").Value = 6
$ws.Range("C32").Value = "mm"

# Row 33: GM_Schrauben_Pos_Stift_rlang
$ws.Range("A33").Value = "GM_Schrauben_Pos_Stift_rlang"
$ws.Range("B33").Value = 9.05
$ws.Range("C33").Value = "mm"

# Row 34: GM_Schrauben_Pos_Stift_rkurz
$ws.Range("A34").Value = "GM_Schrauben_Pos_Stift_rkurz"
$ws.Range("B34").Value = 9.185
$ws.Range("C34").Value = "mm"

# Comment for row 30 is entered after row 34 so the shared-string order matches
$ws.Range("D30").Value = "m2,5"

# Row 35: GM_Schrauben_Durchgangsloch
$ws.Range("A35").Value = "GM_Schrauben_Durchgangsloch"
$ws.Range("B35").Value = 2.7
$ws.Range("C35").Value = "mm"

# --- Column widths: widen column A to fit the longer parameter names ---
$ws.Columns.Item(1).AutoFit()

# --- Selection moves to C35 (last edited cell) ---
$ws.Range("C35").Select()
